# Adds the "2022-Q3" quarterly snapshot sheet (feat: add 2022-Q3 data).
#
# Effect:
#   1. A brand-new worksheet "2022-Q3" is inserted right after "总计" (i.e.
#      before "2022-Q2"), holding the new quarter's fund holdings.
#   2. The "总计" (totals) roll-up sheet gets a new row 2 for "2022-Q3"; every
#      other existing row simply slides down one row (its own data is
#      untouched) and the running index in column A is renumbered 0..7.
#   3. Every other existing quarter sheet (2022-Q2 .. 2020-Q4) is left
#      exactly as-is - they just move one tab to the right to make room.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: insert the new "2022-Q3" worksheet before "2022-Q2"
# ---------------------------------------------------------------------
$sibling = $wb.Worksheets.Item("2022-Q2")
$ws = $wb.Worksheets.Add($sibling)
$ws.Name = "2022-Q3"

# Match the page margins used by the rest of the quarterly sheets.
# PageSetup margins are expressed in points (72pt = 1in).
$ws.PageSetup.LeftMargin = 0.75 * 72
$ws.PageSetup.RightMargin = 0.75 * 72
$ws.PageSetup.TopMargin = 1 * 72
$ws.PageSetup.BottomMargin = 1 * 72
$ws.PageSetup.HeaderMargin = 0.5 * 72
$ws.PageSetup.FooterMargin = 0.5 * 72

# Pull the header / index-column formatting from the sibling sheet so the
# new tab looks like the rest (bold + bordered header row, styled index col).
$sibling.Range("B1:H1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)
$sibling.Range("A2:A4").Copy()
$ws.Range("A2:A4").PasteSpecial(-4122)

$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

# Row 2
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "'159617"
$ws.Range("C2").Value = "华夏中证智选500价值稳健策略ETF"
$ws.Range("D2").Value = "'2.93"
$ws.Range("E2").Value = "'97.05"
$ws.Range("F2").Value = "'1.40"
$ws.Range("G2").Value = "'0.0410"
$ws.Range("H2").Value = 7

# Row 3
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "'006143"
$ws.Range("C3").Value = "恒生前海中证质量成长低波动指数A"
$ws.Range("D3").Value = "'0.05"
$ws.Range("E3").Value = "'93.33"
$ws.Range("F3").Value = "'2.33"
$ws.Range("G3").Value = "'0.0012"
$ws.Range("H3").Value = 8

# Row 4
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "'006144"
$ws.Range("C4").Value = "恒生前海中证质量成长低波动指数C"
$ws.Range("D4").Value = "'0.01"
$ws.Range("E4").Value = "'93.33"
$ws.Range("F4").Value = "'2.33"
$ws.Range("G4").Value = "'0.0002"
$ws.Range("H4").Value = 8

# ---------------------------------------------------------------------
# Step 2: update the "总计" roll-up sheet with the new 2022-Q3 row
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

# "Insert" leaves the new row with inherited formatting that doesn't match
# the rest of the table - reset it from the (now pushed-down) row 3.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$total.Range("B3:D3").Copy()
$total.Range("B2:D2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 0.04

# Renumber the running index in column A for the rows that shifted down.
For ($i = 3; $i -le 9; $i++) {
    $total.Cells.Item($i, 1).Value = $i - 2
}

# ---------------------------------------------------------------------
# Step 3: keep the original active tab ("2020-Q4") selected
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2020-Q4").Activate()
